$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (rows 2-21), columns A,B,C
$data = @(
    @(1.136592060327531, -0.5052947402000403, 1.006910741329194),
    @(-3.102525949478149, -7.741808414459229, 0.0443248748779296),
    @(-0.8462124168872824, -1.410848729312418, -6.25395178794861),
    @(4.942803740501414, -7.510328069329272, -1.564183712005606),
    @(3.684623420238498, -6.547261834144596, 2.809367418289182),
    @(-2.489120721817007, -3.182153344154361, 2.89736366271973),
    @(-2.860559403896316, -2.504158109426502, 1.881697505712511),
    @(2.841274738311764, -2.992074728012082, 1.633560657501217),
    @(0.4830425977706894, -0.9799425601959225, -0.1168481409549758),
    @(2.629005432128906, -1.44602632522583, 6.930900573730469),
    @(0.08325040340423584, -1.20862039923668, 1.943817764520645),
    @(-1.441726684570314, -0.9328206181526182, 5.954500854015357),
    @(-0.886786460876462, 0.562802791595463, 2.521360546350466),
    @(0.6777331829071001, -1.399440765380847, 0.4171198606491046),
    @(1.204387292265894, -3.619132399559024, -0.313209235668175),
    @(1.067639499902725, -3.350790739059448, -3.138641357421875),
    @(3.527233093976973, -5.226221919059754, -3.07070302963257),
    @(7.073627471923828, -3.28519868850708, -0.3319654464721679),
    @(2.874565660953519, -0.3619521260261513, -1.012392401695252),
    @(-6.492176651954649, -0.7134745419025422, 7.389736890792848)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove the now-unused last row (22) that previously held data
$ws.Range("A22:C22").ClearContents()
